$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-DateText {
    param($addr, $text)
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $text
    $c.Style = "Normal"
}

# Row 3: 28/07/2022 -> 28-07-2022; D 0->1; G 0->1
Set-DateText "A3" "28-07-2022"
$ws.Range("D3").Value = 1
$ws.Range("G3").Value = 1

# Row 4: 01/08/2022 -> 01-08-2022; D 0->1; E 0->1; H 1->0
Set-DateText "A4" "01-08-2022"
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 1
$ws.Range("H4").Value = 0

# Row 5: 04/08/2022 -> 04-08-2022
Set-DateText "A5" "04-08-2022"

# Row 6: 08/08/2022 -> 08-08-2022
Set-DateText "A6" "08-08-2022"

# Row 7: 11/08/2022 -> 11-08-2022; D 0->2; E 0->1; F 0->1; H 1->0
Set-DateText "A7" "11-08-2022"
$ws.Range("D7").Value = 2
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 1
$ws.Range("H7").Value = 0

# Row 8: 15/08/2022 -> 15-08-2022
Set-DateText "A8" "15-08-2022"

# Row 9: 18/08/2022 -> 18-08-2022
Set-DateText "A9" "18-08-2022"

# Row 10: 22/08/2022 -> 22-08-2022
Set-DateText "A10" "22-08-2022"

# Row 11: 25/08/2022 -> 25-08-2022
Set-DateText "A11" "25-08-2022"

# Row 12: 29/08/2022 -> 29-08-2022
Set-DateText "A12" "29-08-2022"

# Row 13: 01/09/2022 -> 01-09-2022; D 0->1; E 0->1; H 1->0
Set-DateText "A13" "01-09-2022"
$ws.Range("D13").Value = 1
$ws.Range("E13").Value = 1
$ws.Range("H13").Value = 0

# Row 14: 05/09/2022 -> 05-09-2022
Set-DateText "A14" "05-09-2022"

# Row 15: 08/09/2022 -> 08-09-2022
Set-DateText "A15" "08-09-2022"

# Row 16: 12/09/2022 -> 12-09-2022
Set-DateText "A16" "12-09-2022"

# Row 17: 15/09/2022 -> 15-09-2022
Set-DateText "A17" "15-09-2022"

# Row 18: 19/09/2022 -> 19-09-2022
Set-DateText "A18" "19-09-2022"

# Row 19: 22/09/2022 -> 22-09-2022
Set-DateText "A19" "22-09-2022"

# Row 20: 26/09/2022 -> 26-09-2022
Set-DateText "A20" "26-09-2022"

# Row 21: 29/09/2022 -> 29-09-2022
Set-DateText "A21" "29-09-2022"
